$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($rng, [string]$val) {
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

# Row 2
$ws.Range("D2").Value = "43.080.64"
$ws.Range("E2").Value = "  +0.40%  "

# Row 3
$ws.Range("D3").Value = "2.367.57"
$ws.Range("E3").Value = "  +1.26%  "

# Row 5
Set-TextValue $ws.Range("D5") "302.76"
$ws.Range("E5").Value = "  -0.25%  "

# Row 6
Set-TextValue $ws.Range("D6") "95.21"
$ws.Range("E6").Value = "  +1.08%  "

# Row 7
$ws.Range("E7").Value = "  -0.03%  "

# Row 9
Set-TextValue $ws.Range("D9") "0.480"
$ws.Range("E9").Value = "  -2.85%  "

# Row 10
Set-TextValue $ws.Range("D10") "34.28"
$ws.Range("E10").Value = "  +0.63%  "

# Row 11
$ws.Range("E11").Value = "  +3.15%  "

# Row 12
Set-TextValue $ws.Range("D12") "0.0786"
$ws.Range("E12").Value = "  +0.52%  "

# Row 13
Set-TextValue $ws.Range("D13") "18.23"
$ws.Range("E13").Value = "  -2.62%  "

# Row 14
Set-TextValue $ws.Range("D14") "6.74"
$ws.Range("E14").Value = "  +0.52%  "

# Row 15
$ws.Range("D15").Value = "2.737.43"
$ws.Range("E15").Value = "  +1.25%  "

# Row 16
$ws.Range("D16").Value = "2.359.71"
$ws.Range("E16").Value = "  +1.09%  "

# Row 17
Set-TextValue $ws.Range("D17") "0.796"
$ws.Range("E17").Value = "  +0.41%  "

# Row 18
$ws.Range("D18").Value = "43.135.94"
$ws.Range("E18").Value = "  +0.66%  "

# Row 19
Set-TextValue $ws.Range("D19") "11.94"
$ws.Range("E19").Value = "  -0.99%  "

# Row 20
Set-TextValue $ws.Range("D20") "6.25"
$ws.Range("E20").Value = "  +0.59%  "

# Row 21
$ws.Range("D21").Value = "0.0₃0886"
$ws.Range("E21").Value = "  -0.14%  "

# Row 22
Set-TextValue $ws.Range("D22") "67.91"
$ws.Range("E22").Value = "  -0.01%  "

# Row 23
Set-TextValue $ws.Range("D23") "235.46"
$ws.Range("E23").Value = "  +0.03%  "

# Row 24
$ws.Range("E24").Value = "  -1.18%  "

# Row 25
$ws.Range("E25").Value = "  +0.28%  "

# Row 26
$ws.Range("E26").Value = "  +0.00%  "

# Row 27
Set-TextValue $ws.Range("D27") "24.41"
$ws.Range("E27").Value = "  -0.54%  "

# Row 28
$ws.Range("E28").Value = "  +15.21%  "

# Row 29
Set-TextValue $ws.Range("D29") "9.34"
$ws.Range("E29").Value = "  +2.49%  "

# Row 30
Set-TextValue $ws.Range("D30") "32.00"
$ws.Range("E30").Value = "  +2.30%  "

# Row 31
$ws.Range("E31").Value = "  +0.00%  "

# Row 32
Set-TextValue $ws.Range("D32") "5.01"
$ws.Range("E32").Value = "  +0.72%  "

# Row 33
Set-TextValue $ws.Range("D33") "17.57"
$ws.Range("E33").Value = "  +2.04%  "

# Row 34
Set-TextValue $ws.Range("D34") "0.109"
$ws.Range("E34").Value = "  +7.97%  "

# Row 35
Set-TextValue $ws.Range("D35") "0.0727"
$ws.Range("E35").Value = "  -4.59%  "

# Row 36
$ws.Range("B36").Value = "Monero"
$ws.Range("C36").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
Set-TextValue $ws.Range("D36") "127.00"
$ws.Range("E36").Value = "  +0.91%  "

# Row 37
$ws.Range("B37").Value = "ARBITRUM"
$ws.Range("C37").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
Set-TextValue $ws.Range("D37") "1.83"

# Row 38
Set-TextValue $ws.Range("D38") "2.85"
$ws.Range("E38").Value = "  +3.21%  "

# Row 39
Set-TextValue $ws.Range("D39") "4.30"
$ws.Range("E39").Value = "  -1.53%  "

# Row 40
Set-TextValue $ws.Range("D40") "2.27"
$ws.Range("E40").Value = "  -2.39%  "

# Row 42
Set-TextValue $ws.Range("D42") "20.78"
$ws.Range("E42").Value = "  -5.96%  "

# Row 43
$ws.Range("D43").Value = "1.928.02"
$ws.Range("E43").Value = "  -0.28%  "

# Row 44
$ws.Range("E44").Value = "  -1.19%  "

# Row 45
$ws.Range("E45").Value = "  +2.53%  "

# Row 46
$ws.Range("B46").Value = "FraxShare"
$ws.Range("C46").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
Set-TextValue $ws.Range("D46") "9.21"
$ws.Range("E46").Value = "  -9.45%  "

# Row 47
$ws.Range("B47").Value = "NEARProtocol"
$ws.Range("C47").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
Set-TextValue $ws.Range("D47") "2.72"
$ws.Range("E47").Value = "  +0.64%  "

# Row 48
$ws.Range("D48").Value = "2.598.04"
$ws.Range("E48").Value = "  +1.22%  "

# Row 49
$ws.Range("E49").Value = "  +2.28%  "

# Row 50
Set-TextValue $ws.Range("D50") "71.36"
$ws.Range("E50").Value = "  -0.29%  "

# Row 51
Set-TextValue $ws.Range("D51") "51.30"
$ws.Range("E51").Value = "  -2.77%  "
